# Finished Week 13 logging
# Update the target depth data for row "H" (A2) on both the OFF and DEF sheets.

$wb = $excel.ActiveWorkbook

# OFF sheet (Offense)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 238
$wsOff.Range("C2").Value = 165
$wsOff.Range("D2").Value = 28
$wsOff.Range("E2").Value = 14

# DEF sheet (Defense)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 282
$wsDef.Range("C2").Value = 181
$wsDef.Range("D2").Value = 77
$wsDef.Range("E2").Value = 36
$wsDef.Range("G2").Value = 3
